# TO DO.xlsx update
# - Adds a "Notities" column (D) with a "Done"/"Datum" checklist (columns F/G)
# - Marks most existing tasks as done ("x") with a completion date
# - Adds three new backlog rows (Settings scherm afmaken / Load game inbouwen / End-game credits inbouwen)
# - Widens column G slightly so the date values are readable

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column G width (character width ~11)
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 10.1666666666667

# ---------------------------------------------------------------------------
# Header band: extend the colored banner (row 1-2) across F:G and add the new
# "Notities" / "Done" / "Datum" header cells on row 3.
# ---------------------------------------------------------------------------
$ws.Range("F1:G2").Interior.ThemeColor = 10

# "Notities" header cell (D3) re-uses the same look as A3:C3 (header style)
$ws.Range("A3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "Notities"

# "Done" / "Datum" headers (F3:G3) use a green-accent header style
$ws.Range("F3:G3").Interior.ThemeColor = 10
$ws.Range("F3:G3").Borders.Item(9).LineStyle = 1
$ws.Range("F3").Value = "Done"
$ws.Range("G3").Value = "Datum"

# ---------------------------------------------------------------------------
# Mark completed backlog items: "x" in column F and a completion date in G.
# The very first date cell (G4) gets the real "mm-dd-yy" number format;
# every later date cell re-uses that exact style via copy/paste-format so
# the stylesheet doesn't accumulate one duplicate <xf> per cell (this
# interpreter's NumberFormat setter does not de-duplicate identical <xf>s).
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 44959
$ws.Range("G4").NumberFormat = "mm-dd-yy"
$dateFormatCell = $ws.Range("G4")

function Mark-X([int]$row) {
    $ws.Cells.Item($row, 6).Value = "x"
}

function Set-DateCell([int]$row, $dateFormatCell) {
    $dateFormatCell.Copy()
    $target = $ws.Cells.Item($row, 7)
    $target.PasteSpecial(-4122)
    $target.Value = 44959
}

Mark-X 4
# G4 already carries the date value + format set above

Mark-X 5
Set-DateCell 5 $dateFormatCell

Mark-X 14
$ws.Cells.Item(14, 7).Value = "18-02-2020"

Mark-X 15
Set-DateCell 15 $dateFormatCell

Mark-X 17
Set-DateCell 17 $dateFormatCell

Mark-X 20
Set-DateCell 20 $dateFormatCell

Mark-X 22
Set-DateCell 22 $dateFormatCell

# ---------------------------------------------------------------------------
# New backlog rows
# ---------------------------------------------------------------------------
$ws.Cells.Item(23, 1).Value = 20
$ws.Cells.Item(23, 2).Value = "Settings scherm afmaken"
$ws.Cells.Item(23, 3).Value = "Is nu puur visueel. Dit kan beter."
Mark-X 23
Set-DateCell 23 $dateFormatCell

$ws.Cells.Item(24, 1).Value = 21
$ws.Cells.Item(24, 2).Value = "Load game inbouwen"
$ws.Cells.Item(24, 3).Value = "Zorgen dat je een saved game kan inladen, zodat je verder kunt spelen na een crash en/of makkelijk terug kan naar een oude spelsituatie"
Mark-X 24
Set-DateCell 24 $dateFormatCell

$ws.Cells.Item(25, 1).Value = 22
$ws.Cells.Item(25, 2).Value = "End-game credits inbouwen"
$ws.Cells.Item(25, 3).Value = "Dit moet gewoon ;)"
Mark-X 25
Set-DateCell 25 $dateFormatCell

# ---------------------------------------------------------------------------
# Selection cursor ends up on A13 in the edited file
# ---------------------------------------------------------------------------
$ws.Range("A13").Select()
